# Apply the daily crypto-price refresh to Sheet1.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# Some Price values (column D) look like plain numbers to Excel (e.g. "1.00",
# "14.80", "0.0330"); left alone, COM auto-converts such text to a real number
# and mangles formatting (trailing zeros, scientific notation, etc). To keep
# them as literal text we temporarily force Text format, assign the value, then
# restore the cell to the default "Normal" style so no extra formatting lingers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.785.11"
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "3.110.19"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.65%  "

$ws.Range("E8").Value = "  +7.29%  "

$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("D10").Value = "3.106.92"
$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.729"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").Value = "91.548.03"
$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.70%  "

$ws.Range("D17").Value = "3.674.23"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "3.102.24"
$ws.Range("E18").Value = "  -0.24%  "

$ws.Range("E19").Value = "  -3.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "

# Row 22 now holds a different coin (rows reordered upstream)
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "445.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.32%  "

# Row 23 now holds a different coin (rows reordered upstream)
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.86%  "

# Row 26 now holds a different coin (rows reordered upstream)
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.60%  "

# Row 27 now holds a different coin (rows reordered upstream)
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").Value = "3.264.06"
$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = "  +12.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.227"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.65%  "

$ws.Range("E32").Value = "  -6.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.166"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("E38").Value = "  -5.55%  "

$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "482.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.433"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.86"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.698"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0330"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
